# Applies the cryptos.xlsx price/volume/coin-name update described by the
# commit "Updated cryptos list on Wed Nov 22 17:33:13 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells are stored as plain text in the source file (e.g.
# "36.577.05" uses dots as thousands separators, not a decimal point), so the
# whole Price column is forced to text entry first - otherwise values such as
# "231.83" or "5.16" would silently be re-interpreted as numbers. The format is
# cleared again afterwards so the cells keep the same (default) style as before.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = "36.577.05"
$ws.Range("E2").Value = "  -2.04%  "
$ws.Range("D3").Value = "2.037.75"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "231.83"
$ws.Range("E5").Value = "  -9.45%  "
$ws.Range("D6").Value = "0.600"
$ws.Range("E6").Value = "  -1.42%  "
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("D8").Value = "55.35"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  -1.86%  "
$ws.Range("D10").Value = "56.95"
$ws.Range("E10").Value = "  +0.98%  "
$ws.Range("E11").Value = "  -1.33%  "
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").Value = "2.336.78"
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").Value = "14.32"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("E15").Value = "  -7.67%  "
$ws.Range("D16").Value = "0.761"
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("D17").Value = "5.16"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "2.040.78"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "36.729.40"
$ws.Range("E19").Value = "  -1.35%  "
$ws.Range("D20").Value = "5.85"
$ws.Range("E20").Value = "  +15.77%  "
$ws.Range("D21").Value = "67.64"
$ws.Range("E21").Value = "  -2.98%  "
$ws.Range("D22").Value = "0.0₃0798"
$ws.Range("E22").Value = "  -3.60%  "
$ws.Range("D23").Value = "220.68"
$ws.Range("E23").Value = "  -5.67%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("E26").Value = "  -5.68%  "
$ws.Range("D27").Value = "162.81"
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "18.94"
$ws.Range("E29").Value = "  -2.35%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "0.126"
$ws.Range("E30").Value = "  -5.46%  "
$ws.Range("E31").Value = "  +3.45%  "
$ws.Range("D32").Value = "0.116"
$ws.Range("E32").Value = "  -1.47%  "
$ws.Range("E33").Value = "  -4.01%  "
$ws.Range("E34").Value = "  -2.91%  "
$ws.Range("E35").Value = "  +4.15%  "
$ws.Range("E36").Value = "  -1.79%  "
$ws.Range("E38").Value = "  -2.87%  "
$ws.Range("D39").Value = "5.78"
$ws.Range("E39").Value = "  +8.38%  "
$ws.Range("D40").Value = "3.21"
$ws.Range("E40").Value = "  -5.65%  "
$ws.Range("B41").Value = "HuobiToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D41").Value = "2.95"
$ws.Range("E41").Value = "  -3.64%  "
$ws.Range("B42").Value = "FTXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D42").Value = "4.46"
$ws.Range("E42").Value = "  +33.96%  "
$ws.Range("D43").Value = "1.474.37"
$ws.Range("E43").Value = "  +2.36%  "
$ws.Range("D44").Value = "0.0944"
$ws.Range("E44").Value = "  +3.33%  "
$ws.Range("D45").Value = "93.67"
$ws.Range("E45").Value = "  +5.07%  "
$ws.Range("E46").Value = "  -1.50%  "
$ws.Range("D47").Value = "1.11"
$ws.Range("E47").Value = "  -4.64%  "
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("E49").Value = "  -2.15%  "
$ws.Range("D50").Value = "2.90"
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("E51").Value = "  +2.25%  "

$priceCol.ClearFormats()
